# Major system update: Import Excel data, fix ingredient prices,
# add order statuses, increase quantities x10
#
# rolls.xlsx - sheet1: add a "cost" column (D), populate sale_price (C)
# and cost (D) for every product row, re-order a handful of product
# names in the tail of the table, and drop the now-unused last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Remove the trailing row (id 24 "кляр для темпура") ---------------
$ws.Rows(25).Delete()

# --- 2) Fix the rotated product names for rows 20-24 ----------------------
$ws.Range("B20").Value = "чедр ролл"
$ws.Range("B21").Value = "мини рол огурец"
$ws.Range("B22").Value = "соус сушиза"
$ws.Range("B23").Value = "спайси соус"
$ws.Range("B24").Value = "сырный соус"

# --- 3) Add the "cost" header in D1, matching the style of C1 -------------
$ws.Range("D1").Value = "cost"
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)  # xlPasteFormats

# --- 4) Populate sale_price (C) and cost (D) for every product row --------
$prices = @(
    @(96, 94.59),
    @(203, 203.21),
    @(129, 129.13),
    @(79, 81.37),
    @(139, 141.9),
    @(160, 163.29),
    @(82, 84.42),
    @(141, 143.12),
    @(196, 196.52),
    @(209, 209.71),
    @(82, 82.25),
    @(203, 203.11),
    @(74, 76.38),
    @(100, 100.7),
    @(45, 45.85),
    @(86, 86.29000000000001),
    @(48, 48.72),
    @(61, 140.2),
    @(118, 118.17),
    @(34, 311.65),
    @(0, 524.9),
    @(0, 352),
    @(0, 348.5)
)

$row = 2
foreach ($pair in $prices) {
    $ws.Cells.Item($row, 3).Value = $pair[0]
    $ws.Cells.Item($row, 4).Value = $pair[1]
    $row = $row + 1
}

Write-Host "applied rolls.xlsx update"
